$d = $word.ActiveDocument

# 1. Reviewer name change: "Джабраилов Х. А." -> "Фурлетов Ю. М."
$d.Content.Find.Execute("Джабраилов Х. А.", $true, $false, $false, $false, $false, $true, 1, $false, "Фурлетов Ю. М.", 2) | Out-Null

# 2. Merge the split "202" / "3" runs into a single "2023" run.
$d.Content.Find.Execute("2023", $true, $false, $false, $false, $false, $true, 1, $false, "2023", 2) | Out-Null

# 3. Locate the paragraph that now holds "2023" followed by the page break.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "2023*") {
        $targetIdx = $i
        break
    }
}

$yearPara = $d.Paragraphs.Item($targetIdx)

# 4. Insert a brand-new paragraph right after it for the repository link.
$yearPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($targetIdx + 1)

# Paragraph formatting to match the target layout.
$newPara.Format.Alignment = 3          # wdAlignParagraphJustify ("both")
$newPara.Format.LineSpacingRule = 0    # wdLineSpaceSingle/auto w/ explicit value below
$newPara.Format.LineSpacing = 18       # 18 pt -> 360 twips
$newPara.Format.FirstLineIndent = 42.55 # 42.55 pt -> 851 twips

# 5. Fill the paragraph with label + URL text, then convert the URL into a
#    real hyperlink (Hyperlinks.Add always anchors to the start of the
#    paragraph, so build the plain text first and target only the URL
#    substring).
$label = "Ссылка на репозиторий: "
$url = "https://github.com/PatriotRossii/informatics_laboratories_8"

$nr = $newPara.Range
$nr.Text = $label + $url

$pStart = $newPara.Range.Start
$pEnd = $newPara.Range.End
$urlStart = $pStart + $label.Length
$urlEnd = $pEnd - 1

$urlRange = $d.Range($urlStart, $urlEnd)
$d.Hyperlinks.Add($urlRange, $url) | Out-Null

Write-Output "done"
